# Update "想去人数" (wanted-to-go headcount) values in column F
# on the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row = newValue }
$updates = @{
    "展览" = @{
        2  = 171
        6  = 478
        8  = 168
        9  = 2460
        10 = 148
        12 = 167
        16 = 320
        24 = 138
        25 = 42
        26 = 1543
        27 = 21
        28 = 385
        29 = 392
        30 = 194
        32 = 392
    }
    "全部类型" = @{
        2  = 171
        7  = 478
        9  = 168
        10 = 2460
        11 = 148
        13 = 167
        17 = 320
        25 = 138
        26 = 42
        27 = 1543
        28 = 21
        29 = 385
        30 = 392
        31 = 194
        33 = 392
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
